# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh updates to H:N columns
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 86
$ws.Range("H86").Value = 5160.6
$ws.Range("I86").Value = 5200.75
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 5200.75
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -4077.75
$ws.Range("N86").Value = -7246
# Row 89
$ws.Range("H89").Value = 5160.6
$ws.Range("I89").Value = 5200.75
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 26003.75
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -20387.75
$ws.Range("N89").Value = -36232
# Row 111
$ws.Range("H111").Value = 1116
$ws.Range("I111").Value = 800
$ws.Range("J111").Value = 1326.6666
$ws.Range("K111").Value = 2400
$ws.Range("L111").Value = 3979.9998
$ws.Range("M111").Value = 667
$ws.Range("N111").Value = -10113.9998
# Row 132
$ws.Range("H132").Value = 1211520.6
$ws.Range("I132").Value = 1253846
$ws.Range("J132").Value = 5250
$ws.Range("K132").Value = 3761538
$ws.Range("L132").Value = 15750
$ws.Range("M132").Value = -3759008
$ws.Range("N132").Value = -20810
# Row 138
$ws.Range("H138").Value = 16398168
$ws.Range("I138").Value = 38463016
$ws.Range("J138").Value = 7138.8
$ws.Range("K138").Value = 115389048
$ws.Range("L138").Value = 21416.4
$ws.Range("M138").Value = -115383908
$ws.Range("N138").Value = -31696.4

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6385.318
$ws.Range("I32").Value = 3674.2388
$ws.Range("J32").Value = 16476.555
$ws.Range("K32").Value = 3674.2388
$ws.Range("L32").Value = 16476.555
$ws.Range("M32").Value = -3387.2388
$ws.Range("N32").Value = -17050.555
# Row 74
$ws.Range("H74").Value = 815.3488
$ws.Range("I74").Value = 730.2593000000001
$ws.Range("J74").Value = 958.9375
$ws.Range("K74").Value = 730.2593000000001
$ws.Range("L74").Value = 958.9375
$ws.Range("M74").Value = 143.7406999999999
$ws.Range("N74").Value = -2706.9375
# Row 77
$ws.Range("H77").Value = 815.3488
$ws.Range("I77").Value = 730.2593000000001
$ws.Range("J77").Value = 958.9375
$ws.Range("K77").Value = 3651.2965
$ws.Range("L77").Value = 4794.6875
$ws.Range("M77").Value = 716.7034999999996
$ws.Range("N77").Value = -13530.6875
# Row 122
$ws.Range("H122").Value = 2194.2683
$ws.Range("I122").Value = 1495.9524
$ws.Range("J122").Value = 2927.5
$ws.Range("K122").Value = 4487.857199999999
$ws.Range("L122").Value = 8782.5
$ws.Range("M122").Value = -2037.857199999999
$ws.Range("N122").Value = -13682.5
# Row 134
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 1536.0638
$ws.Range("I134").Value = 1245.9667
$ws.Range("J134").Value = 2048
$ws.Range("K134").Value = 3737.9001
$ws.Range("L134").Value = 6144
$ws.Range("M134").Value = -1202.9001
$ws.Range("N134").Value = -11214

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 4902.1333
$ws.Range("I16").Value = 2491.5
$ws.Range("J16").Value = 7657.143
$ws.Range("K16").Value = 2491.5
$ws.Range("L16").Value = 7657.143
$ws.Range("M16").Value = -2204.5
$ws.Range("N16").Value = -8231.143
# Row 99
$ws.Range("H99").Value = 3388.889
$ws.Range("I99").Value = 2500
$ws.Range("K99").Value = 2500
$ws.Range("M99").Value = -1002
# Row 113
$ws.Range("H113").Value = 4902.1333
$ws.Range("I113").Value = 2491.5
$ws.Range("J113").Value = 7657.143
$ws.Range("K113").Value = 2491.5
$ws.Range("L113").Value = 7657.143
$ws.Range("M113").Value = -321.5
$ws.Range("N113").Value = -11997.143
# Row 122
$ws.Range("H122").Value = 1553.3334
$ws.Range("I122").Value = 1483.3334
$ws.Range("J122").Value = 1833.3334
$ws.Range("K122").Value = 4450.0002
$ws.Range("L122").Value = 5500.0002
$ws.Range("M122").Value = -2000.0002
$ws.Range("N122").Value = -10400.0002
# Row 126
$ws.Range("H126").Value = 3388.889
$ws.Range("I126").Value = 2500
$ws.Range("K126").Value = 7500
$ws.Range("M126").Value = -5030

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 137
$ws.Range("H137").Value = 34668610
$ws.Range("J137").Value = 2224270.2
$ws.Range("L137").Value = 6672810.600000001
$ws.Range("N137").Value = -6683010.600000001
# Row 138
$ws.Range("H138").Value = 10871971
$ws.Range("I138").Value = 11907082
$ws.Range("K138").Value = 35721246
$ws.Range("M138").Value = -35716106

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2953.6
$ws.Range("I102").Value = 1560
$ws.Range("J102").Value = 4727.273
$ws.Range("K102").Value = 1560
$ws.Range("L102").Value = 4727.273
$ws.Range("M102").Value = 62
$ws.Range("N102").Value = -7971.273
# Row 122
$ws.Range("H122").Value = 2091.5454
$ws.Range("I122").Value = 1858.1428
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 5574.428400000001
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -3124.428400000001
$ws.Range("N122").Value = -12400

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 309.88
$ws.Range("I22").Value = 308.6111
$ws.Range("J22").Value = 313.14285
$ws.Range("K22").Value = 308.6111
$ws.Range("L22").Value = 313.14285
$ws.Range("M22").Value = -13.61110000000002
$ws.Range("N22").Value = -903.14285
# Row 27
$ws.Range("H27").Value = 309.88
$ws.Range("I27").Value = 308.6111
$ws.Range("J27").Value = 313.14285
$ws.Range("K27").Value = 308.6111
$ws.Range("L27").Value = 313.14285
$ws.Range("M27").Value = -201.6111
$ws.Range("N27").Value = -527.14285
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()
# Row 61
$ws.Range("H61").Value = 1774.6666
$ws.Range("I61").Value = 1767
$ws.Range("J61").Value = 1777.2222
$ws.Range("K61").Value = 1767
$ws.Range("L61").Value = 1777.2222
$ws.Range("M61").Value = -1565
$ws.Range("N61").Value = -2181.2222
# Row 113
$ws.Range("H113").Value = 1774.6666
$ws.Range("I113").Value = 1767
$ws.Range("J113").Value = 1777.2222
$ws.Range("K113").Value = 1767
$ws.Range("L113").Value = 1777.2222
$ws.Range("M113").Value = 403
$ws.Range("N113").Value = -6117.2222
# Row 115
$ws.Range("H115").Value = 10302
$ws.Range("J115").Value = 10302
$ws.Range("L115").Value = 10302
$ws.Range("N115").Value = -12652
# Row 118
$ws.Range("H118").Value = 75995
$ws.Range("J118").Value = 75995
$ws.Range("L118").Value = 75995
$ws.Range("N118").Value = -79309
# Row 121
$ws.Range("H121").Value = 34183.227
$ws.Range("J121").Value = 34183.227
$ws.Range("L121").Value = 34183.227
$ws.Range("N121").Value = -37677.227
# Row 122
$ws.Range("H122").Value = 3162.2703
$ws.Range("I122").Value = 3000.121
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 9000.363000000001
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -6550.363000000001
$ws.Range("N122").Value = -18400
# Row 123
$ws.Range("H123").Value = 19072
$ws.Range("J123").Value = 19072
$ws.Range("L123").Value = 19072
$ws.Range("N123").Value = -28872
# Row 125
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
# Row 127
$ws.Range("H127").Value = 35000
$ws.Range("J127").Value = 35000
$ws.Range("L127").Value = 35000
$ws.Range("N127").Value = -44920

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 6146.778
$ws.Range("I107").Value = 331.66666
$ws.Range("J107").Value = 9054.333000000001
$ws.Range("K107").Value = 994.9999799999999
$ws.Range("L107").Value = 27162.999
$ws.Range("M107").Value = 925.0000200000001
$ws.Range("N107").Value = -31002.999
# Row 122
$ws.Range("H122").Value = 1968.3334
$ws.Range("I122").Value = 1950
$ws.Range("J122").Value = 2005
$ws.Range("K122").Value = 5850
$ws.Range("L122").Value = 6015
$ws.Range("M122").Value = -3400
$ws.Range("N122").Value = -10915
